# Update the "Dictionary" sheet's Short name (column C) descriptions with
# fuller, more descriptive endpoint names, per commit "added new endpoint
# descriptions".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dictionary")

$ws.Range("C2").Value  = "Any effect except mortality"
$ws.Range("C3").Value  = "Any effect in 5 days"
$ws.Range("C4").Value  = "Any effect in 24 hours"
$ws.Range("C5").Value  = "5 day total movement"
$ws.Range("C9").Value  = "Delayed development at 24h"
$ws.Range("C14").Value = "5 day behavior transition"
$ws.Range("C16").Value = "Notochord at 5 days"
$ws.Range("C19").Value = "Spontaneous motion at 24h"

# Column C needs to widen to fit the new, longer text (matches the
# resulting OOXML's stored column width of 27).
$ws.Columns("C").ColumnWidth = 26.166666666666668

# The sheet was left active with C19 (the last-edited cell) selected.
$ws.Activate()
$ws.Range("C19").Select()
